$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-5 from 2023-10-13 to 2023-10-22
$ws.Range("C2:C5").Value = 45221
